# obj_ins_statistics.pptx - apply the authored changes:
#  1. Merge the split "Object" / "-level" runs in the "Object-level human masks"
#     caption back into a single "Object-level" run.
#  2. Add five new Times New Roman, 16pt caption textboxes (Conversation, Din,
#     Singing, Playing Instruments, Monologue) on top of the picture.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- 1. Fix the "Object" + "-level" split run -----------------------------
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
        $full = $shp.TextFrame.TextRange.Text
        if ($full -eq "Object-level human masks") {
            $shp.TextFrame.TextRange.Characters(1, 12).Text = "Object-level"
        }
    }
}

# --- 2. Add the five new caption textboxes ---------------------------------
# AddTextbox() hands out the lowest free shape id. The real deck this change
# came from already had ids 14-16 in use elsewhere, so the five new labels
# landed on ids 17-21; reproduce that numbering here by first "burning"
# through the currently-free low ids (2,3,4,14,15,16) with throwaway
# textboxes, then deleting them before creating the real labels.
$burn = @()
for ($k = 0; $k -lt 6; $k++) {
    $burn += $s.Shapes.AddTextbox(1, 0, 0, 10, 10)
}
foreach ($d in $burn) {
    $d.Delete()
}

function New-Caption($name, $left, $top, $width, $height, $text) {
    $tb = $s.Shapes.AddTextbox(1, $left, $top, $width, $height)
    $tb.Name = $name
    $tb.Fill.Visible = $false
    $tb.TextFrame.WordWrap = $true
    $tb.TextFrame.AutoSize = 1
    $tr = $tb.TextFrame.TextRange
    $tr.Text = $text
    $tr.Font.Size = 16
    $tr.Font.Name = "Times New Roman"
    $tr.Font.NameFarEast = "Times New Roman"
    $tr.Font.NameComplexScript = "Times New Roman"
    return $tb
}

New-Caption "TextBox 16" 166.49212598425197 23.99291338582677 103.788031496063 26.65779527559055 "Conversation" | Out-Null
New-Caption "TextBox 17" 502.97330708661417 119.53173228346456 46.99314960629921 26.65779527559055 "Din" | Out-Null
New-Caption "TextBox 18" 339.85141732283466 276.7872440944882 68.49275590551181 26.65779527559055 "Singing" | Out-Null
New-Caption "TextBox 19" 602.800157480315 338.0074803149606 144.83181102362204 26.65779527559055 "Playing Instruments" | Out-Null
New-Caption "TextBox 20" 801.0307086614173 217.8855905511811 91.788031496063 26.65779527559055 "Monologue" | Out-Null
